$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update ILAsimLIN (row 16): new D16 value; B16 = D16/D2 recalculates automatically
$ws.Range("D16").Value = 10830.3666970787

# Update OLGsolveLIN (row 20): formula now references C20/C2 instead of D20/D2
$ws.Range("C20").Value = 0.039299027062952498
$ws.Range("B20").Formula = "=C20/C2"

# Update OLGsimLIN (row 23): formula now references C23/C2 instead of D23/D2
$ws.Range("C23").Value = 2177.3693303610198
$ws.Range("B23").Formula = "=C23/C2"

# Update the selected cell in the bottom-right (scrollable) pane of the frozen view
$ws.Activate()
$ws.Range("B16").Select()
